{"js": "// \"Empty key bug fixed\": the paragraph held a placeholder/garbled string;\n// replace it with the real text and add a manual line break (w:br) right\n// after it, matching the author's committed change.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\n// Replace the existing run's text in place (keeps the run's rPr, i.e. the\n// ru-RU language tag, on the text run).\nparagraph.getRange(\"Content\").insertText(\"\u0432\u044c\u0432\u0445 \u044f\u0448\u043a \u0435\u0448\u0436\u0448\u044c \u0446\u044f\u0448\u043a\", Word.InsertLocation.replace);\n\n// Append a manual line break as a new run at the end of the paragraph.\nparagraph.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\nawait context.sync();\n\n// The new break run is created without character formatting; stamp it with\n// the same ru-RU language tag the original run had so it matches\n// <w:rPr><w:lang w:val=\"ru-RU\"/></w:rPr> on the <w:br/> run.\nconst endRange = paragraph.getRange(\"End\");\nendRange.languageId = \"ru-RU\";\nawait context.sync();\n", "ps1": "# \"Empty key bug fixed\": the paragraph held a placeholder/garbled string;\n# replace it with the real text and add a manual line break (w:br) right\n# after it, matching the author's committed change.\n$d = $word.ActiveDocument\n\n$p = $d.Paragraphs(1)\n\n# Replace the existing run's text in place (keeps the run's rPr, i.e. the\n# ru-RU language tag, on the text run). Paragraph.Range.Text includes the\n# trailing paragraph mark, so assigning Text replaces the visible content\n# while Word keeps the paragraph mark.\n$p.Range.Text = \"\u0432\u044c\u0432\u0445 \u044f\u0448\u043a \u0435\u0448\u0436\u0448\u044c \u0446\u044f\u0448\u043a\"\n\n# Re-fetch the paragraph (its Range end shifted after the text replace) and\n# collapse a range right before the paragraph mark, then insert a manual\n# line break (wdLineBreak = 6) there -> new run with <w:br/>.\n$p = $d.Paragraphs(1)\n$insertionPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)\n$insertionPoint.InsertBreak(6)\n\n# The new break run is created without character formatting; stamp it with\n# the same ru-RU language tag the original run had so it matches\n# <w:rPr><w:lang w:val=\"ru-RU\"/></w:rPr> on the <w:br/> run.\n$p = $d.Paragraphs(1)\n$endRange = $d.Range($p.Range.End - 1, $p.Range.End - 1)\n$endRange.LanguageID = \"ru-RU\"\n"}
